# Updates cryptos list: prices, volume(1h) percentages, and two row swaps
# (row44/row45 Monero<->EnergySwap, row46/row47 LidoDAOToken<->ARBITRUM)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @("D2", "66.063.68"),
    @("E2", "  -1.05%  "),
    @("D3", "3.745.76"),
    @("E3", "  +0.57%  "),
    @("E4", "  +0.16%  "),
    @("D5", "405.35"),
    @("E5", "  -4.36%  "),
    @("D6", "129.13"),
    @("E6", "  -1.96%  "),
    @("D7", "3.733.12"),
    @("E7", "  +0.47%  "),
    @("D8", "0.603"),
    @("E8", "  -5.83%  "),
    @("E9", "  +0.03%  "),
    @("D10", "0.720"),
    @("E10", "  -5.88%  "),
    @("E11", "  -8.86%  "),
    @("D12", "0.0000355"),
    @("E12", "  -6.46%  "),
    @("D13", "40.35"),
    @("E13", "  -5.55%  "),
    @("D14", "4.333.86"),
    @("E14", "  +0.65%  "),
    @("D15", "9.62"),
    @("E15", "  -6.49%  "),
    @("D16", "14.40"),
    @("E16", "  +9.75%  "),
    @("E17", "  -1.75%  "),
    @("D18", "3.739.46"),
    @("E18", "  -0.06%  "),
    @("D19", "19.39"),
    @("E19", "  -7.18%  "),
    @("D20", "66.354.00"),
    @("E20", "  -0.63%  "),
    @("E21", "  -6.63%  "),
    @("D22", "407.12"),
    @("E22", "  -8.92%  "),
    @("D23", "14.29"),
    @("E23", "  -9.94%  "),
    @("D24", "85.00"),
    @("E24", "  -5.95%  "),
    @("D25", "3.01"),
    @("E25", "  -5.40%  "),
    @("D26", "35.93"),
    @("E26", "  -5.85%  "),
    @("D27", "5.62"),
    @("E27", "  +13.15%  "),
    @("D28", "3.09"),
    @("E28", "  -6.47%  "),
    @("D29", "9.30"),
    @("E29", "  -8.77%  "),
    @("D30", "12.32"),
    @("E30", "  -2.30%  "),
    @("E31", "  -2.70%  "),
    @("E32", "  -3.88%  "),
    @("D33", "7.10"),
    @("E33", "  -2.20%  "),
    @("D34", "0.154"),
    @("E34", "  -5.17%  "),
    @("D35", "38.80"),
    @("E35", "  -7.58%  "),
    @("D36", "0.999"),
    @("E36", "  -0.05%  "),
    @("D37", "55.14"),
    @("E37", "  -2.40%  "),
    @("D38", "0.0₃0725"),
    @("E38", "  -0.86%  "),
    @("E39", "  -7.63%  "),
    @("D40", "2.86"),
    @("E40", "  -6.03%  "),
    @("D41", "0.998"),
    @("E41", "  +0.23%  "),
    @("E42", "  -8.61%  "),
    @("D43", "3.17"),
    @("E43", "  +20.83%  "),
    @("B44", "EnergySwap"),
    @("C44", "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"),
    @("D44", "26.84"),
    @("E44", "  -8.41%  "),
    @("B45", "Monero"),
    @("C45", "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"),
    @("D45", "144.92"),
    @("E45", "  -0.68%  "),
    @("B46", "LidoDAOToken"),
    @("C46", "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"),
    @("D46", "3.23"),
    @("E46", "  -6.94%  "),
    @("B47", "ARBITRUM"),
    @("C47", "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"),
    @("D47", "2.04"),
    @("E47", "  -4.86%  "),
    @("D48", "2.81"),
    @("E48", "  -3.57%  "),
    @("E49", "  -3.70%  "),
    @("E50", "  -4.64%  "),
    @("E51", "  -6.07%  ")
)

foreach ($pair in $changes) {
    $cell = $ws.Range($pair[0])
    # Force text format so numeric-looking strings (e.g. "405.35", "0.720")
    # keep their exact original text representation instead of being
    # auto-converted to numbers by Excel.
    $cell.NumberFormat = "@"
    $cell.Value = $pair[1]
}
